$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tokens")

# New "Latex Rule" column (G) header
$ws.Range("G1").Value = "Latex Rule"

# New Latex-escaped grammar rule for the CELL token (row 9), forced to text
# via a leading apostrophe so Excel stores it with the quote-prefix style
# (mirrors the other regex-like cells in this sheet, e.g. C2/E8/C19/E13).
$ws.Range("G9").Value = "'[a-zA-Z0-9]+ "":"" ([0-9A-Z\_.]+ `$\mid`$ ""'"" ([0-9A-Z\_ !@\#\`$\%\textasciicircum{}\&*()\-+=\{\}:;`$\mid`$\textless\textgreater,./?\textbackslash{}\textbackslash{}] `$\mid`$ ``''')+  ``\ '\ ') ``!'"

# Leave the selection where the author left it after entering the new data
$ws.Range("G10").Select()
